# timesheet.xlsx edit: fill in the second punch (in/out) for the
# Sunday row of week 1 (row 13) and add task notes for the last two
# rows of that week (rows 13 and 14).
#
# Row 13 (Sun:) originally only had a morning in/out (C13/D13). A second
# in/out pair is recorded for the evening: in at 10:20 PM, out at 11:15 PM.
# Row 14 (Mon:) gets a note describing the work performed; row 13 also
# gets a note. The Pay Hours / weekly / monthly totals are formula-driven
# and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second "In"/"Out" punch for row 13 (values are Excel day-fractions,
# matching the existing time-only number format already applied to the
# F:G columns of that row).
$ws.Range("E13").Value = 0.93055555555555547   # 10:20:00 PM
$ws.Range("F13").Value = 0.96875               # 11:15:00 PM

# Notes column (I) for the last two rows of week 1.
$ws.Range("I13").Value = "Grabbing computer and creating window installation media"
$ws.Range("I14").Value = "Installing windows, creating account, downloading data, and getting instructions for task"

# Leave the active cell on E14, matching where the user's cursor ended up.
[void]$ws.Range("E14").Select()
